{"js": "// Locate the M2Doc field \" m:''.sampleTable() \" (rendered as a Word complex\n// field: fldChar begin / instrText... / fldChar end) and rewrite it as\n// plain literal text runs \"{m:''.sampleTable()}\" (TokenIteratorFieldRewriterSplit\n// turns the field into a leading \"{\" run, the instruction-text runs turned\n// into plain text runs, and a trailing \"}\" run), keeping the orange run\n// coloring on \"''.\", \"sample\" and \"Table()\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet targetParagraph = null;\nlet targetField = null;\n\nfor (const paragraph of paragraphs.items) {\n  const fields = paragraph.fields;\n  fields.load(\"items\");\n  await context.sync();\n\n  for (const field of fields.items) {\n    field.load(\"code\");\n  }\n  await context.sync();\n\n  for (const field of fields.items) {\n    if (field.code.trim() === \"m:''.sampleTable()\") {\n      targetParagraph = paragraph;\n      targetField = field;\n      break;\n    }\n  }\n  if (targetField) {\n    break;\n  }\n}\n\nif (targetField) {\n  // Removes every run belonging to the field (fldChar begin/end +\n  // instrText runs) in one shot, leaving the host paragraph empty.\n  targetField.delete();\n  await context.sync();\n\n  const colorRpr =\n    '<w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr>';\n  const innerRuns =\n    \"<w:r><w:t>{</w:t></w:r>\" +\n    \"<w:r><w:t>m</w:t></w:r>\" +\n    \"<w:r><w:t>:</w:t></w:r>\" +\n    \"<w:r>\" + colorRpr + \"<w:t>''.</w:t></w:r>\" +\n    \"<w:r>\" + colorRpr + \"<w:t>sample</w:t></w:r>\" +\n    \"<w:r>\" + colorRpr + \"<w:t>Table()</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>';\n\n  const packageXml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" + innerRuns + \"</w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\";\n\n  const targetRange = targetParagraph.getRange(\"Whole\");\n  targetRange.insertOoxml(packageXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the M2Doc field \" m:''.sampleTable() \" (rendered as a Word complex\n# field: fldChar begin / instrText.../ fldChar end) and rewrite it as plain\n# literal text runs \"{m:''.sampleTable()}\" (TokenIteratorFieldRewriterSplit\n# splits the field into a leading \"{\" run, the instruction-text runs turned\n# into plain text runs, and a trailing \"}\" run), keeping the orange run\n# coloring on \"''.\", \"sample\" and \"Table()\".\n\n$apos = [char]39\n\n$found = $false\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n    $field = $d.Fields.Item($i)\n    if ($field.Code.Text.Trim() -eq \"m:''.sampleTable()\") {\n        $found = $true\n\n        # Find the 1-based Document.Paragraphs index that hosts the field\n        # (Range.Paragraphs is not reliably scoped in this host, so locate\n        # it by comparing character offsets instead).\n        $codeStart = $field.Code.Start\n        $fieldParagraphIndex = -1\n        for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {\n            $cand = $d.Paragraphs.Item($j)\n            if ($codeStart -ge $cand.Range.Start -and $codeStart -lt $cand.Range.End) {\n                $fieldParagraphIndex = $j\n                break\n            }\n        }\n\n        # Field.Delete() removes every run belonging to the field\n        # (fldChar begin/end + instrText runs) in one shot.\n        $field.Delete()\n\n        $p = $d.Paragraphs.Item($fieldParagraphIndex)\n        $insertRng = $p.Range\n        # Exclude the trailing paragraph mark from the replacement range.\n        $insertRng.End = $insertRng.End - 1\n\n        $xmlFrag = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n            '<w:r><w:t>{</w:t></w:r>' + `\n            '<w:r><w:t>m</w:t></w:r>' + `\n            '<w:r><w:t>:</w:t></w:r>' + `\n            '<w:r><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>' + $apos + $apos + '.</w:t></w:r>' + `\n            '<w:r><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>sample</w:t></w:r>' + `\n            '<w:r><w:rPr><w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/></w:rPr><w:t>Table()</w:t></w:r>' + `\n            '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' + `\n            '</w:p>'\n\n        $null = $insertRng.InsertXML($xmlFrag)\n        break\n    }\n}\n\nif (-not $found) {\n    Write-Output \"WARNING: sampleTable field not found\"\n}\n"}
